$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coureurs")

# Update the Australia score column (H) for drivers that scored points
$ws.Range("H2").Value = 25
$ws.Range("H3").Value = 11
$ws.Range("H7").Value = 18
$ws.Range("H10").Value = 4
$ws.Range("H11").Value = 8
$ws.Range("H13").Value = 2
$ws.Range("H14").Value = 12
$ws.Range("H15").Value = 15
$ws.Range("H17").Value = 6
$ws.Range("H19").Value = 1

# Update the active selection on this sheet
$ws.Activate()
$ws.Range("H24").Select()
